# [FIX] update data formatting for consistency
#
# The "income" sheet's header row (row 1) labels several balance-sheet
# columns ambiguously. Clarify that these are aggregate/total figures,
# and simplify the combined liabilities+equity header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Total Current Assets"
$ws.Range("E1").Value = "Total Liabilities"
$ws.Range("F1").Value = "Total Equity"
$ws.Range("L1").Value = "Liabilities and Equity"
